# Rename stock symbol "APPL" to "AAPL" (typo fix) across all worksheets.
# APPL was always meant to be AAPL.

$wb = $excel.ActiveWorkbook

$ws_rsu = $wb.Worksheets.Item("rsu")
$ws_rsu.Range("B6").Value = "AAPL"

$ws_dividends = $wb.Worksheets.Item("dividends")
$ws_dividends.Range("B3").Value = "AAPL"

$ws_sell = $wb.Worksheets.Item("sell_orders")
$ws_sell.Range("B6").Value = "AAPL"
$ws_sell.Range("B7").Value = "AAPL"

$ws_rsu.Activate()
